$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 154
$ws.Range("I6").Value = 185
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 555
$ws.Range("L6").Value = 90
$ws.Range("M6").Value = -443
$ws.Range("N6").Value = -314

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 254
$ws.Range("I39").Value = 254
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 762
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -466

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3779.5557
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 4002
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 4002
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -4352

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2233.1667
$ws.Range("I132").Value = 2233.1667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6699.500100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4169.500100000001

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7146175.5
$ws.Range("I138").Value = 33337000
$ws.Range("J138").Value = 3223.182
$ws.Range("K138").Value = 100011000
$ws.Range("L138").Value = 9669.545999999998
$ws.Range("M138").Value = -100005860
$ws.Range("N138").Value = -19949.546

# ARM row 16
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 1199
$ws.Range("I16").Value = 1199
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1199
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -912
$ws.Range("N16").ClearContents()

# ARM row 26
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 2883.3333
$ws.Range("I26").Value = 1925
$ws.Range("J26").Value = 4800
$ws.Range("K26").Value = 1925
$ws.Range("L26").Value = 4800
$ws.Range("M26").Value = -1595
$ws.Range("N26").Value = -5460

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2864
$ws.Range("I122").Value = 2864
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8592
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6142
$ws.Range("N122").ClearContents()

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3115.75
$ws.Range("I132").Value = 2671.7273
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 8015.1819
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -5485.1819
$ws.Range("N132").Value = -29060

# BSM row 12
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 837
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 996.25
$ws.Range("K12").Value = 200
$ws.Range("L12").Value = 996.25
$ws.Range("M12").Value = -32
$ws.Range("N12").Value = -1332.25

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1445.7273
$ws.Range("I94").Value = 1499.3
$ws.Range("J94").Value = 910
$ws.Range("K94").Value = 1499.3
$ws.Range("L94").Value = 910
$ws.Range("M94").Value = -1048.3
$ws.Range("N94").Value = -1812

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1116.125
$ws.Range("I99").Value = 1148.7142
$ws.Range("J99").Value = 888
$ws.Range("K99").Value = 1148.7142
$ws.Range("L99").Value = 888
$ws.Range("M99").Value = 349.2858000000001
$ws.Range("N99").Value = -3884

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2470
$ws.Range("I105").Value = 2455
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 2455
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -708
$ws.Range("N105").Value = -5994

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2434.2727
$ws.Range("I134").Value = 2177.7
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6533.099999999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -3998.099999999999
$ws.Range("N134").Value = -20070

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3638.05
$ws.Range("I31").Value = 1884.1333
$ws.Range("J31").Value = 8899.799999999999
$ws.Range("K31").Value = 1884.1333
$ws.Range("L31").Value = 8899.799999999999
$ws.Range("M31").Value = -1589.1333
$ws.Range("N31").Value = -9489.799999999999

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3638.05
$ws.Range("I34").Value = 1884.1333
$ws.Range("J34").Value = 8899.799999999999
$ws.Range("K34").Value = 1884.1333
$ws.Range("L34").Value = 8899.799999999999
$ws.Range("M34").Value = -1682.1333
$ws.Range("N34").Value = -9303.799999999999

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2499
$ws.Range("I62").Value = 2499
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2499
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1875
$ws.Range("N62").ClearContents()

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2499
$ws.Range("I65").Value = 2499
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 12495
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -9375
$ws.Range("N65").ClearContents()

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2774.6
$ws.Range("I105").Value = 2468.25
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 2468.25
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -721.25
$ws.Range("N105").Value = -7494

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9334.666999999999
$ws.Range("I132").Value = 5802.2
$ws.Range("J132").Value = 13750.25
$ws.Range("K132").Value = 17406.6
$ws.Range("L132").Value = 41250.75
$ws.Range("M132").Value = -14876.6
$ws.Range("N132").Value = -46310.75

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1998.5454
$ws.Range("I134").Value = 1580.6666
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 4741.9998
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -2206.9998
$ws.Range("N134").Value = -12570

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 100044.1
$ws.Range("I33").Value = 48.875
$ws.Range("J33").Value = 500025
$ws.Range("K33").Value = 293.25
$ws.Range("L33").Value = 3000150
$ws.Range("M33").Value = -10.25
$ws.Range("N33").Value = -3000716

# CUL row 38
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 497.84616
$ws.Range("I38").Value = 186.88889
$ws.Range("J38").Value = 1197.5
$ws.Range("K38").Value = 560.6666700000001
$ws.Range("L38").Value = 3592.5
$ws.Range("M38").Value = -213.6666700000001
$ws.Range("N38").Value = -4286.5

# CUL row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3320.2727
$ws.Range("I117").Value = 5000
$ws.Range("J117").Value = 3152.3
$ws.Range("K117").Value = 15000
$ws.Range("L117").Value = 9456.900000000001
$ws.Range("M117").Value = -11558
$ws.Range("N117").Value = -16340.9

# GSM row 26
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 30000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 30000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 30000
$ws.Range("N26").Value = -30560

# GSM row 50
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H50").Value = 30000
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 30000
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -30996

# GSM row 53
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 23000
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 23000
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 23000
$ws.Range("N53").Value = -24262

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2795.8
$ws.Range("I80").Value = 1993
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 1993
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -995
$ws.Range("N80").Value = -5996

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2795.8
$ws.Range("I83").Value = 1993
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 9965
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -4973
$ws.Range("N83").Value = -29984

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3417.4546
$ws.Range("I102").Value = 1755.625
$ws.Range("J102").Value = 7849
$ws.Range("K102").Value = 1755.625
$ws.Range("L102").Value = 7849
$ws.Range("M102").Value = -133.625
$ws.Range("N102").Value = -11093

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1000
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -530

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2963.182
$ws.Range("I132").Value = 2963.182
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8889.545999999998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6359.545999999998

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 7318.6
$ws.Range("I136").Value = 5534.3335
$ws.Range("J136").Value = 9995
$ws.Range("K136").Value = 16603.0005
$ws.Range("L136").Value = 29985
$ws.Range("M136").Value = -14053.0005
$ws.Range("N136").Value = -35085

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 700
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2100
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 430

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8990.333000000001
$ws.Range("I136").Value = 8990.333000000001
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 26970.999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -24420.999
